# Daily attendance processing - 2026-01-04 17:01:21
# Swap the "Recorded By" (column G) value order from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# for every row in the active worksheet where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
